# Actualización automática 2025-09-26 11:50:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO sales for MEGAMAFERS S.A. (row 13) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M13").Value = 10280.02

# --- Sheet "VENTA MENSUAL": septiembre sales for MEGAMAFERS S.A. (row 13) + total (row 23) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F13").Value = 11639.25
$wsMensual.Range("F23").Value = 51113.75

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO group (row 12) and TOTAL (row 15) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D12").Value = 43904.91
$wsCumplimiento.Range("E12").Value = -7081.266907882906
$wsCumplimiento.Range("F12").Value = 1.192302181784909

$wsCumplimiento.Range("D15").Value = 50670.61
$wsCumplimiento.Range("E15").Value = 4754.13316613377
$wsCumplimiento.Range("F15").Value = 0.9142236319998196
